$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1428.3846
$ws.Range("I2").Value = 1616.2727
$ws.Range("J2").Value = 395
$ws.Range("K2").Value = 1616.2727
$ws.Range("L2").Value = 395
$ws.Range("M2").Value = -1503.2727
$ws.Range("N2").Value = -621

$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

$ws.Range("H33").Value = 1164.909
$ws.Range("I33").Value = 1201.375
$ws.Range("K33").Value = 1201.375
$ws.Range("M33").Value = -972.375

$ws.Range("H55").Value = 129.05882
$ws.Range("I55").Value = 82.666664
$ws.Range("J55").Value = 181.25
$ws.Range("K55").Value = 82.666664
$ws.Range("L55").Value = 181.25
$ws.Range("M55").Value = 131.333336
$ws.Range("N55").Value = -609.25

$ws.Range("H64").Value = 4225.4688
$ws.Range("I64").Value = 4075.3076
$ws.Range("J64").Value = 4328.2104
$ws.Range("K64").Value = 4075.3076
$ws.Range("L64").Value = 4328.2104
$ws.Range("M64").Value = -3827.3076
$ws.Range("N64").Value = -4824.2104

$ws.Range("H67").Value = 4225.4688
$ws.Range("I67").Value = 4075.3076
$ws.Range("J67").Value = 4328.2104
$ws.Range("K67").Value = 4075.3076
$ws.Range("L67").Value = 4328.2104
$ws.Range("M67").Value = -3217.3076
$ws.Range("N67").Value = -6044.2104

$ws.Range("H132").Value = 3841.6072
$ws.Range("I132").Value = 2617.4285
$ws.Range("J132").Value = 7514.143
$ws.Range("K132").Value = 7852.2855
$ws.Range("L132").Value = 22542.429
$ws.Range("M132").Value = -5322.2855
$ws.Range("N132").Value = -27602.429

$ws.Range("H138").Value = 2012.9423
$ws.Range("I138").Value = 909.3200000000001
$ws.Range("J138").Value = 3034.8147
$ws.Range("K138").Value = 2727.96
$ws.Range("L138").Value = 9104.444100000001
$ws.Range("M138").Value = 2412.04
$ws.Range("N138").Value = -19384.4441

$ws.Range("H141").Value = 3189.625
$ws.Range("I141").Value = 3074.5715
$ws.Range("J141").Value = 3995
$ws.Range("K141").Value = 9223.7145
$ws.Range("L141").Value = 11985
$ws.Range("M141").Value = -4043.7145
$ws.Range("N141").Value = -22345

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 468.41666
$ws.Range("J5").Value = 823.3333
$ws.Range("L5").Value = 823.3333
$ws.Range("N5").Value = -1047.3333

$ws.Range("H25").Value = 1933.3334
$ws.Range("I25").Value = 900
$ws.Range("K25").Value = 900
$ws.Range("M25").Value = -498

$ws.Range("H132").Value = 2448.1606
$ws.Range("I132").Value = 2288
$ws.Range("J132").Value = 2885.9333
$ws.Range("K132").Value = 6864
$ws.Range("L132").Value = 8657.7999
$ws.Range("M132").Value = -4334
$ws.Range("N132").Value = -13717.7999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 468.41666
$ws.Range("J4").Value = 823.3333
$ws.Range("L4").Value = 823.3333
$ws.Range("N4").Value = -1053.3333

$ws.Range("H134").Value = 2220.2122
$ws.Range("I134").Value = 1158.0454
$ws.Range("J134").Value = 4344.5454
$ws.Range("K134").Value = 3474.1362
$ws.Range("L134").Value = 13033.6362
$ws.Range("M134").Value = -939.1361999999999
$ws.Range("N134").Value = -18103.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 21.3
$ws.Range("I7").Value = 12.2
$ws.Range("J7").Value = 30.4
$ws.Range("K7").Value = 12.2
$ws.Range("L7").Value = 30.4
$ws.Range("M7").Value = 100.8
$ws.Range("N7").Value = -256.4

$ws.Range("H31").Value = 2894.125
$ws.Range("I31").Value = 1611.2727
$ws.Range("J31").Value = 5716.4
$ws.Range("K31").Value = 1611.2727
$ws.Range("L31").Value = 5716.4
$ws.Range("M31").Value = -1316.2727
$ws.Range("N31").Value = -6306.4

$ws.Range("H34").Value = 2894.125
$ws.Range("I34").Value = 1611.2727
$ws.Range("J34").Value = 5716.4
$ws.Range("K34").Value = 1611.2727
$ws.Range("L34").Value = 5716.4
$ws.Range("M34").Value = -1409.2727
$ws.Range("N34").Value = -6120.4

$ws.Range("H132").Value = 2294.375
$ws.Range("I132").Value = 978.7857
$ws.Range("J132").Value = 3317.611
$ws.Range("K132").Value = 2936.3571
$ws.Range("L132").Value = 9952.832999999999
$ws.Range("M132").Value = -406.3571000000002
$ws.Range("N132").Value = -15012.833

$ws.Range("H135").Value = 53268.75
$ws.Range("J135").Value = 34700
$ws.Range("L135").Value = 34700
$ws.Range("N135").Value = -44840

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 62.692307
$ws.Range("I2").Value = 19.333334
$ws.Range("J2").Value = 99.85714
$ws.Range("K2").Value = 116.000004
$ws.Range("L2").Value = 599.14284
$ws.Range("M2").Value = -3.000004000000004
$ws.Range("N2").Value = -825.14284

$ws.Range("H38").Value = 47.88889
$ws.Range("I38").Value = 35.25
$ws.Range("J38").Value = 84
$ws.Range("K38").Value = 105.75
$ws.Range("L38").Value = 252
$ws.Range("M38").Value = 241.25
$ws.Range("N38").Value = -946

$ws.Range("H107").Value = 1062.375
$ws.Range("I107").Value = 300
$ws.Range("J107").Value = 1171.2858
$ws.Range("K107").Value = 900
$ws.Range("L107").Value = 3513.8574
$ws.Range("M107").Value = 1020
$ws.Range("N107").Value = -7353.857400000001

$ws.Range("H118").Value = 3170.923
$ws.Range("I118").Value = 907
$ws.Range("J118").Value = 4177.1113
$ws.Range("K118").Value = 2721
$ws.Range("L118").Value = 12531.3339
$ws.Range("M118").Value = -1478
$ws.Range("N118").Value = -15017.3339

$ws.Range("H122").Value = 26316332
$ws.Range("J122").Value = 1000.4286
$ws.Range("L122").Value = 9003.857399999999
$ws.Range("N122").Value = -13903.8574

$ws.Range("H131").Value = 1354.4746
$ws.Range("J131").Value = 1288.7307
$ws.Range("L131").Value = 3866.1921
$ws.Range("N131").Value = -13946.1921

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 5975
$ws.Range("I45").Value = 5000
$ws.Range("J45").Value = 6170
$ws.Range("K45").Value = 5000
$ws.Range("L45").Value = 6170
$ws.Range("M45").Value = -4441
$ws.Range("N45").Value = -7288

$ws.Range("H51").Value = 29533.334
$ws.Range("J51").Value = 29533.334
$ws.Range("L51").Value = 29533.334
$ws.Range("N51").Value = -30551.334

$ws.Range("H124").Value = 28390
$ws.Range("J124").Value = 28390
$ws.Range("L124").Value = 28390
$ws.Range("N124").Value = -38210

$ws.Range("H133").Value = 27222.25
$ws.Range("I133").Value = 24900
$ws.Range("J133").Value = 27996.334
$ws.Range("K133").Value = 24900
$ws.Range("L133").Value = 27996.334
$ws.Range("M133").Value = -19840
$ws.Range("N133").Value = -38116.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 24252
$ws.Range("I25").Value = 8500
$ws.Range("J25").Value = 40004
$ws.Range("K25").Value = 8500
$ws.Range("L25").Value = 40004
$ws.Range("M25").Value = -8270
$ws.Range("N25").Value = -40464

$ws.Range("H46").Value = 1476.3077
$ws.Range("I46").Value = 2150
$ws.Range("J46").Value = 898.8570999999999
$ws.Range("K46").Value = 2150
$ws.Range("L46").Value = 898.8570999999999
$ws.Range("M46").Value = -1962
$ws.Range("N46").Value = -1274.8571

$ws.Range("H117").Value = 37500
$ws.Range("J117").Value = 37500
$ws.Range("L117").Value = 37500
$ws.Range("N117").Value = -46678

$ws.Range("H118").Value = 28000
$ws.Range("J118").Value = 28000
$ws.Range("L118").Value = 28000
$ws.Range("N118").Value = -31314

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H132").Value = 3154.1143
$ws.Range("I132").Value = 2195.2856
$ws.Range("J132").Value = 4592.357
$ws.Range("K132").Value = 6585.8568
$ws.Range("L132").Value = 13777.071
$ws.Range("M132").Value = -4055.8568
$ws.Range("N132").Value = -18837.071

$ws.Range("H136").Value = 3120.0408
$ws.Range("I136").Value = 1998.1
$ws.Range("J136").Value = 4891.5264
$ws.Range("K136").Value = 5994.299999999999
$ws.Range("L136").Value = 14674.5792
$ws.Range("M136").Value = -3444.299999999999
$ws.Range("N136").Value = -19774.5792

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 39333.332
$ws.Range("J125").Value = 39333.332
$ws.Range("L125").Value = 39333.332
$ws.Range("N125").Value = -49173.332

$ws.Range("H136").Value = 14100579
$ws.Range("I136").Value = 27055782
$ws.Range("J136").Value = 2271.3235
$ws.Range("K136").Value = 81167346
$ws.Range("L136").Value = 6813.970499999999
$ws.Range("M136").Value = -81164796
$ws.Range("N136").Value = -11913.9705
